{"js": "// Update the date heading and the 25 division problems/answers in the\n// worksheet table (5 data rows x 5 columns, each data row followed by\n// three blank spacer rows).\n\n// 1) Update the date paragraph (\"2024-04-28 Sunday\" -> \"2024-04-29 Monday\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst datePara = paragraphs.items[0];\ndatePara.load(\"text\");\nawait context.sync();\n\nif (datePara.text === \"2024-04-28 Sunday\") {\n  datePara.insertText(\"2024-04-29 Monday\", \"Replace\");\n}\n\n// 2) Update the table of division problems.\n// Each entry gives the 0-based table row index (data rows are 0, 4, 8,\n// 12, 16 \u2014 every 4th row, the others are blank spacer rows) and the five\n// new cell values, left to right.\nconst rowUpdates = [\n  { row: 0, values: [\"92\u00f72=46, 0\", \"81\u00f77=11, 4\", \"15\u00f78=1, 7\", \"15\u00f74=3, 3\", \"93\u00f74=23, 1\"] },\n  { row: 4, values: [\"75\u00f78=9, 3\", \"35\u00f72=17, 1\", \"70\u00f77=10, 0\", \"22\u00f73=7, 1\", \"13\u00f78=1, 5\"] },\n  { row: 8, values: [\"44\u00f74=11, 0\", \"78\u00f79=8, 6\", \"83\u00f79=9, 2\", \"62\u00f76=10, 2\", \"43\u00f76=7, 1\"] },\n  { row: 12, values: [\"90\u00f79=10, 0\", \"99\u00f78=12, 3\", \"64\u00f76=10, 4\", \"23\u00f76=3, 5\", \"46\u00f72=23, 0\"] },\n  { row: 16, values: [\"74\u00f72=37, 0\", \"40\u00f77=5, 5\", \"82\u00f73=27, 1\", \"44\u00f75=8, 4\", \"42\u00f75=8, 2\"] },\n];\n\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nfor (const { row, values } of rowUpdates) {\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(row, col);\n    cell.value = values[col];\n  }\n}\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division problems/answers in the\n# worksheet table (5 data rows x 5 columns, each data row followed by\n# three blank spacer rows).\n\n$d = $word.ActiveDocument\n\n# 1) Update the date paragraph (\"2024-04-28 Sunday\" -> \"2024-04-29 Monday\").\n$find = $d.Content.Find\n$find.Text = \"2024-04-28 Sunday\"\n$find.Execute(\"2024-04-28 Sunday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2024-04-29 Monday\", 2)\n\n# 2) Update the table of division problems.\n# Each row gives the 1-based table row index (data rows are 1, 5, 9, 13,\n# 17 \u2014 every 4th row, the others are blank spacer rows) and the five new\n# cell values, left to right.\n$t = $d.Tables.Item(1)\n\n$rowUpdates = @(\n    @{ Row = 1;  Values = @(\"92\u00f72=46, 0\", \"81\u00f77=11, 4\", \"15\u00f78=1, 7\", \"15\u00f74=3, 3\", \"93\u00f74=23, 1\") },\n    @{ Row = 5;  Values = @(\"75\u00f78=9, 3\", \"35\u00f72=17, 1\", \"70\u00f77=10, 0\", \"22\u00f73=7, 1\", \"13\u00f78=1, 5\") },\n    @{ Row = 9;  Values = @(\"44\u00f74=11, 0\", \"78\u00f79=8, 6\", \"83\u00f79=9, 2\", \"62\u00f76=10, 2\", \"43\u00f76=7, 1\") },\n    @{ Row = 13; Values = @(\"90\u00f79=10, 0\", \"99\u00f78=12, 3\", \"64\u00f76=10, 4\", \"23\u00f76=3, 5\", \"46\u00f72=23, 0\") },\n    @{ Row = 17; Values = @(\"74\u00f72=37, 0\", \"40\u00f77=5, 5\", \"82\u00f73=27, 1\", \"44\u00f75=8, 4\", \"42\u00f75=8, 2\") }\n)\n\nforeach ($rowUpdate in $rowUpdates) {\n    $rowIndex = $rowUpdate.Row\n    $values = $rowUpdate.Values\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $values[$col - 1]\n    }\n}\n\nWrite-Output \"done\"\n"}
